$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 14 ("ايفا كيراتين حمام كريم"),
# pushing the existing rows 14-18 down to 16-20. Use row 13 as the style
# template so the new rows inherit the same borders/fonts/number formats
# as the rest of the data table.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(15).Insert()

$ws.Range("A13:Q13").Copy($ws.Range("A14:Q14"))
$ws.Range("A13:Q13").Copy($ws.Range("A15:Q15"))

# Row heights matching the rest of the data rows (the row below the
# inserted pair keeps shifting its auto-fit height as text reflows, so
# pin every affected row explicitly to match the refreshed report).
$ws.Rows.Item(14).RowHeight = 25.5
$ws.Rows.Item(15).RowHeight = 24.75
$ws.Rows.Item(16).RowHeight = 25.5
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 24.75

# New row 14: TERRAMYCIN EYE OINT. 5 GM
$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "TERRAMYCIN EYE OINT. 5 GM"
$ws.Range("H14").Value = "3:0"
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value = "1"
$ws.Range("L14").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"
$ws.Range("N14").Value = "28.00"
$ws.Range("P14").NumberFormat = "@"
$ws.Range("P14").Value = "28.0000"
$ws.Range("P14").NumberFormat = "0.00"
$ws.Range("Q14").Value = "1:0"

# New row 15: TOBRIN 0.3% EYE DROPS 5 ML
$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "TOBRIN 0.3% EYE DROPS 5 ML"
$ws.Range("H15").Value = "2:0"
$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = "1"
$ws.Range("L15").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"
$ws.Range("N15").Value = "23.00"
$ws.Range("P15").NumberFormat = "@"
$ws.Range("P15").Value = "23.0000"
$ws.Range("P15").NumberFormat = "0.00"
$ws.Range("Q15").Value = "1:0"

# Renumber the rows that shifted down (previously 8/9/10, now 10/11/12).
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11
$ws.Range("A18").Value = 12

# Update the grand total to include the two new rows (625.43 + 28.00 + 23.00).
$ws.Range("P19").Value = 676.43
